# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt - Plátano"
# as row 339, shifting the existing rows 339:364 down to 340:365.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(339).Insert()

$ws.Cells.Item(339, 1).Value  = 4
$ws.Cells.Item(339, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(339, 3).Value  = "Los Lagos"
$ws.Cells.Item(339, 4).Value  = 44585
$ws.Cells.Item(339, 5).Value  = 10
$ws.Cells.Item(339, 6).Value  = "Fruta"
$ws.Cells.Item(339, 7).Value  = 100108
$ws.Cells.Item(339, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(339, 9).Value  = 100108006
$ws.Cells.Item(339, 10).Value = "Plátano"
$ws.Cells.Item(339, 11).Value = "Sin especificar"
$ws.Cells.Item(339, 12).Value = "Primera Pintón"
$ws.Cells.Item(339, 13).Value = 800
$ws.Cells.Item(339, 14).Value = 20000
$ws.Cells.Item(339, 15).Value = 20000
$ws.Cells.Item(339, 16).Value = 20000
$ws.Cells.Item(339, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(339, 18).Value = "Ecuador"
$ws.Cells.Item(339, 19).Value = 1000
$ws.Cells.Item(339, 20).Value = 20
